$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "numeros"

$ws.Range("E1").Select()
